$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has ended without an agreement.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision process did not result in any agreement for which movie to show, and thus, no decision was made.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for the screening event.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision to select a movie was not reached, so no movie has been acquired for Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for the screening on Friday.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision process concluded without selecting a movie for Friday, resulting in no decision being made.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday has not been made.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been finalized.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("D11").Value = "Barbie_was_selected, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday's showing.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday could not be made.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision-making meeting did not reach a consensus on which movie to show on Friday. As a result, no movie will be acquired for that day.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision process did not lead to an agreement on a movie for Friday, so there is no decision to report.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: No decision was made regarding Friday's movie.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision to acquire the rights for the movie `"Oppenheimer`" has been successfully made.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The rights to both movies `"Barbie`" and `"Oppenheimer`" have been acquired for the Friday showing.`n"
$ws.Range("D26").Value = "both_movies, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision process did not lead to a selection of a movie for Friday, resulting in no decision being made.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on which movie to show on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision to acquire the rights to `"Barbie`" has been confirmed.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was reached regarding the movie for Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.`n"
$ws.Range("D34").Value = "Barbie_was_selected, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on what movie to play on Friday.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding Friday's movie.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: No decision was made about which movie to show on Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding which movie to play on Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired successfully.`n"
$ws.Range("D41").Value = "both_movies, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie has been chosen for Friday.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to show `"Barbie`" on Friday.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected to be shown on Friday.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision process concluded without agreeing on a movie for Friday, resulting in no decision being made.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision-making process did not result in a selection for Friday's movie.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The committee did not reach a decision on which movie to show on Friday.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be shown on Friday as no agreement was reached.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D51").Value = "both_movies, "
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded successfully. `"Oppenheimer`" will be shown on Friday.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision process has concluded without a definitive choice for Friday's movie, resulting in no decision being made.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has resulted in no consensus.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: It appears there has been no decision reached about what movie to show on Friday.`n"
$ws.Range("D59").Value = "no_decision, "
$ws.Range("C60").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has concluded without a selection.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been recorded with no agreement on what movie to show on Friday.`n"
$ws.Range("D62").Value = "no_decision, "
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday ended without a conclusion.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie will be shown on Friday.`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision regarding Friday's movie has not been made.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for the movie `"Oppenheimer`" will be acquired.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision process concluded without reaching an agreement on a movie for Friday.`n"
$ws.Range("C69").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Barbie`" is selected as the movie to acquire for Friday.`n"
$ws.Range("C70").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday ended without a conclusive selection, and thus no decision was made.`n"
